# Script: apply betexplorer scrape update (Atualizado por script em 31-10-2023 20:45)
# - Rows 101-103 and 105-108: match rows were re-sorted/rotated by the upstream scraper;
#   this rewrites columns F:V (home/away teams, scores, odds, timestamps, url) in place.
#   Columns A:E (Indice/pais/torneio/temporada/data_partida) are unaffected.
# - A brand-new match (row 111, Alessandria vs Atalanta U23) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the re-sorted match rows (F:V only; A:E unchanged) ---
# Row 101: Virtus Verona 0 - 2 Mantova
$ws.Range("F101").Value = 'Virtus Verona'
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 'Mantova'
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = 2.21
$ws.Range("K101").Value = '27/10/2023 02:42'
$ws.Range("L101").Value = 2.72
$ws.Range("M101").Value = '28/10/2023 16:07'
$ws.Range("N101").Value = 2.88
$ws.Range("O101").Value = '27/10/2023 02:42'
$ws.Range("P101").Value = 3.09
$ws.Range("Q101").Value = '28/10/2023 16:07'
$ws.Range("R101").Value = 3.31
$ws.Range("S101").Value = '27/10/2023 02:42'
$ws.Range("T101").Value = 2.72
$ws.Range("U101").Value = '28/10/2023 16:07'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/virtus-verona-mantova/ObXsINS4/'

# Row 102: AlbinoLeffe 1 - 0 Arzignano
$ws.Range("F102").Value = 'AlbinoLeffe'
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 'Arzignano'
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2.33
$ws.Range("K102").Value = '27/10/2023 02:42'
$ws.Range("L102").Value = 3.05
$ws.Range("M102").Value = '28/10/2023 16:08'
$ws.Range("N102").Value = 2.88
$ws.Range("O102").Value = '27/10/2023 02:42'
$ws.Range("P102").Value = 2.75
$ws.Range("Q102").Value = '28/10/2023 16:08'
$ws.Range("R102").Value = 3.07
$ws.Range("S102").Value = '27/10/2023 02:42'
$ws.Range("T102").Value = 2.74
$ws.Range("U102").Value = '28/10/2023 16:08'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/albinoleffe-arzignano/l4c544b5/'

# Row 103: Pro Vercelli 2 - 1 Pro Patria
$ws.Range("F103").Value = 'Pro Vercelli'
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 'Pro Patria'
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 2.02
$ws.Range("K103").Value = '27/10/2023 02:42'
$ws.Range("L103").Value = 2.26
$ws.Range("M103").Value = '28/10/2023 16:07'
$ws.Range("N103").Value = 2.99
$ws.Range("O103").Value = '27/10/2023 02:42'
$ws.Range("P103").Value = 3.08
$ws.Range("Q103").Value = '28/10/2023 16:07'
$ws.Range("R103").Value = 3.66
$ws.Range("S103").Value = '27/10/2023 02:42'
$ws.Range("T103").Value = 3.45
$ws.Range("U103").Value = '28/10/2023 16:07'
$ws.Range("V103").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pro-vercelli-pro-patria/fJuRL5Tu/'

# Row 105: Trento 1 - 1 Alessandria
$ws.Range("F105").Value = 'Trento'
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 'Alessandria'
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.74
$ws.Range("K105").Value = '27/10/2023 02:42'
$ws.Range("L105").Value = 1.9
$ws.Range("M105").Value = '28/10/2023 18:21'
$ws.Range("N105").Value = 3.28
$ws.Range("O105").Value = '27/10/2023 02:42'
$ws.Range("P105").Value = 3.05
$ws.Range("Q105").Value = '28/10/2023 18:21'
$ws.Range("R105").Value = 4.54
$ws.Range("S105").Value = '27/10/2023 02:42'
$ws.Range("T105").Value = 4.92
$ws.Range("U105").Value = '28/10/2023 18:21'
$ws.Range("V105").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/trento-alessandria/joZZJqbh/'

# Row 106: Renate 2 - 2 Novara
$ws.Range("F106").Value = 'Renate'
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 'Novara'
$ws.Range("I106").Value = 2
$ws.Range("J106").Value = 1.84
$ws.Range("K106").Value = '27/10/2023 02:42'
$ws.Range("L106").Value = 1.66
$ws.Range("M106").Value = '28/10/2023 18:23'
$ws.Range("N106").Value = 3.18
$ws.Range("O106").Value = '27/10/2023 02:42'
$ws.Range("P106").Value = 3.85
$ws.Range("Q106").Value = '28/10/2023 18:23'
$ws.Range("R106").Value = 4.12
$ws.Range("S106").Value = '27/10/2023 02:42'
$ws.Range("T106").Value = 5.11
$ws.Range("U106").Value = '28/10/2023 18:23'
$ws.Range("V106").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/renate-novara/Y9vVKPrn/'

# Row 107: Giana Erminio 1 - 1 Pro Sesto
$ws.Range("F107").Value = 'Giana Erminio'
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 'Pro Sesto'
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.15
$ws.Range("K107").Value = '27/10/2023 02:42'
$ws.Range("L107").Value = 2.17
$ws.Range("M107").Value = '28/10/2023 18:29'
$ws.Range("N107").Value = 2.93
$ws.Range("O107").Value = '27/10/2023 02:42'
$ws.Range("P107").Value = 3.16
$ws.Range("Q107").Value = '28/10/2023 18:29'
$ws.Range("R107").Value = 3.38
$ws.Range("S107").Value = '27/10/2023 02:42'
$ws.Range("T107").Value = 3.57
$ws.Range("U107").Value = '28/10/2023 18:29'
$ws.Range("V107").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/giana-erminio-pro-sesto/ruC83ODB/'

# Row 108: Pergolettese 3 - 2 Legnago Salus
$ws.Range("F108").Value = 'Pergolettese'
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = 'Legnago Salus'
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 2.27
$ws.Range("K108").Value = '27/10/2023 02:42'
$ws.Range("L108").Value = 2.23
$ws.Range("M108").Value = '27/10/2023 13:42'
$ws.Range("N108").Value = 2.87
$ws.Range("O108").Value = '27/10/2023 02:42'
$ws.Range("P108").Value = 3.11
$ws.Range("Q108").Value = '28/10/2023 16:36'
$ws.Range("R108").Value = 3.2
$ws.Range("S108").Value = '27/10/2023 02:42'
$ws.Range("T108").Value = 3.42
$ws.Range("U108").Value = '27/10/2023 13:42'
$ws.Range("V108").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/pergolettese-legnago-salus/Ic4L0McU/'

# --- Append the new match row 111 ---
# Copy row 110 formatting first (gives A111 the "Indice" style and E111 the date style)
$ws.Range("A110:V110").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A111").Value = 110
$ws.Range("B111").Value = 'italy'
$ws.Range("C111").Value = 'serie-c-group-a'
$ws.Range("D111").Value = '2023-2024'
$ws.Range("E111").Value = 45230.86458333334
$ws.Range("F111").Value = 'Alessandria'
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 'Atalanta U23'
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2.85
$ws.Range("K111").Value = '31/10/2023 08:12'
$ws.Range("L111").Value = 3.5
$ws.Range("M111").Value = '31/10/2023 20:10'
$ws.Range("N111").Value = 2.94
$ws.Range("O111").Value = '31/10/2023 08:12'
$ws.Range("P111").Value = 2.98
$ws.Range("Q111").Value = '31/10/2023 20:10'
$ws.Range("R111").Value = 2.61
$ws.Range("S111").Value = '31/10/2023 08:12'
$ws.Range("T111").Value = 2.29
$ws.Range("U111").Value = '31/10/2023 20:10'
$ws.Range("V111").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-a/alessandria-atalanta/nwdBf8Dq/'

"Done: rows 101-103, 105-108 rewritten; row 111 appended."
